$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update claim numbers in column B
# B2 keeps the same value "1120194100385" (unchanged)
# B3 changes from "1120170200908" to "1120170200917 " (note trailing space)
# Leading apostrophe forces Excel to store the numeric-looking string as text
# (quote-prefixed), matching the original cell formatting.
$ws.Range("B3").Value = "'1120170200917 "

# Update Importe in column C
$ws.Range("C3").Value = 100

# Update the active selection to I9
$ws.Range("I9").Select()
